$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '30.079.02'
$c.Style = "Normal"

$ws.Range('E2').Value = '  -1.51%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.105.67'
$c.Style = "Normal"

$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  -0.74%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '349.85'
$c.Style = "Normal"

$ws.Range('E6').Value = '  -0.66%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5161'
$c.Style = "Normal"

$ws.Range('E7').Value = '  -1.46%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.4470'
$c.Style = "Normal"

$ws.Range('E8').Value = '  -1.38%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '52.63'
$c.Style = "Normal"

$ws.Range('E9').Value = '  -5.25%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.08952'
$c.Style = "Normal"

$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  +0.75%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '25.81'
$c.Style = "Normal"

$ws.Range('E12').Value = '  +5.08%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.105.91'
$c.Style = "Normal"

$ws.Range('E13').Value = '  -0.44%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '8.249'
$c.Style = "Normal"

$ws.Range('E14').Value = '  +1.75%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '6.750'
$c.Style = "Normal"

$ws.Range('E15').Value = '  -1.18%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '99.18'
$c.Style = "Normal"

$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('E18').Value = '  -0.79%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '20.80'
$c.Style = "Normal"

$ws.Range('E19').Value = '  +7.69%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.06664'
$c.Style = "Normal"

$ws.Range('E20').Value = '  -0.37%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.272'
$c.Style = "Normal"

$ws.Range('E22').Value = '  +0.67%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '30.188.31'
$c.Style = "Normal"

$ws.Range('E23').Value = '  -1.35%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '12.91'
$c.Style = "Normal"

$ws.Range('E24').Value = '  +0.77%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.355'
$c.Style = "Normal"

$ws.Range('E25').Value = '  -0.28%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.356.61'
$c.Style = "Normal"

$ws.Range('E26').Value = '  -0.20%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '22.07'
$c.Style = "Normal"

$ws.Range('E27').Value = '  -1.25%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.561'
$c.Style = "Normal"

$ws.Range('E28').Value = '  +2.20%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '162.72'
$c.Style = "Normal"

$ws.Range('E29').Value = '  -0.45%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '133.50'
$c.Style = "Normal"

$ws.Range('E30').Value = '  +0.04%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.180'
$c.Style = "Normal"

$ws.Range('E31').Value = '  -2.68%  '
$ws.Range('E32').Value = '  +0.23%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.643'
$c.Style = "Normal"

$ws.Range('E33').Value = '  +0.69%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '6.277'
$c.Style = "Normal"

$ws.Range('E34').Value = '  -0.94%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '3.975'
$c.Style = "Normal"

$ws.Range('E35').Value = '  +0.30%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.927'
$c.Style = "Normal"

$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('E37').Value = '  -1.79%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.02584'
$c.Style = "Normal"

$ws.Range('E38').Value = '  -0.97%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.06858'
$c.Style = "Normal"

$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '12.83'
$c.Style = "Normal"

$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.2315'
$c.Style = "Normal"

$ws.Range('E41').Value = '  +0.31%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.6841'
$c.Style = "Normal"

$ws.Range('E42').Value = '  +0.23%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.256'
$c.Style = "Normal"

$ws.Range('E43').Value = '  +0.16%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '14.37'
$c.Style = "Normal"

$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.324'
$c.Style = "Normal"

$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.6422'
$c.Style = "Normal"

$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('E47').Value = '  +3.69%  '
$ws.Range('E48').Value = '  -0.16%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '83.60'
$c.Style = "Normal"

$ws.Range('E49').Value = '  +0.64%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.224'
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.07246'
$c.Style = "Normal"

$ws.Range('E51').Value = '  +0.81%  '
